# Apply crypto price/volume updates per the commit diff (Tue Mar 28 08:43:24 UTC 2023 refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.009.86'
$ws.Range('E2').Value = '  -3.06%  '
$ws.Range('D3').Value = '1.726.74'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('D5').Value = '''310.69'
$ws.Range('E5').Value = '  -5.17%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = '''0.4810'
$ws.Range('E7').Value = '  +3.38%  '
$ws.Range('D8').Value = '''0.3485'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '''43.30'
$ws.Range('E9').Value = '  +3.16%  '
$ws.Range('D10').Value = '''0.07242'
$ws.Range('E10').Value = '  -1.57%  '
$ws.Range('D11').Value = '''1.052'
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').Value = '''19.94'
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('D14').Value = '''5.877'
$ws.Range('D15').Value = '1.726.13'
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('D16').Value = '''6.851'
$ws.Range('E16').Value = '  -4.13%  '
$ws.Range('D17').Value = '''86.87'
$ws.Range('E17').Value = '  -5.56%  '
$ws.Range('D18').Value = '''0.00001033'
$ws.Range('D19').Value = '''0.06398'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = '''16.60'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').Value = '''5.707'
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').Value = '27.064.14'
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').Value = '''10.94'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').Value = '''2.072'
$ws.Range('E25').Value = '  -3.59%  '
$ws.Range('D26').Value = '''154.05'
$ws.Range('E26').Value = '  -4.65%  '
$ws.Range('D27').Value = '''19.99'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').Value = '1.921.88'
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('D29').Value = '''2.072'
$ws.Range('E29').Value = '  -3.66%  '
$ws.Range('D30').Value = '''120.78'
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('D31').Value = '''1.044'
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').Value = '''0.09306'
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').Value = '''5.392'
$ws.Range('E34').Value = '  -2.78%  '
$ws.Range('D35').Value = '''0.05926'
$ws.Range('E35').Value = '  -2.46%  '
$ws.Range('D36').Value = '''0.02184'
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('D37').Value = '''1.427'
$ws.Range('E37').Value = '  +5.81%  '
$ws.Range('E38').Value = '  -5.64%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = '''0.1995'
$ws.Range('E39').Value = '  -3.32%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = '''4.765'
$ws.Range('E40').Value = '  -2.62%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').Value = '''0.5983'
$ws.Range('E42').Value = '  -2.55%  '
$ws.Range('D43').Value = '''1.099'
$ws.Range('E43').Value = '  -6.64%  '
$ws.Range('D44').Value = '''7.539'
$ws.Range('E44').Value = '  -2.86%  '
$ws.Range('D45').Value = '''12.77'
$ws.Range('E45').Value = '  -1.78%  '
$ws.Range('D46').Value = '''3.581'
$ws.Range('E46').Value = '  -4.05%  '
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('D48').Value = '''118.69'
$ws.Range('E48').Value = '  -3.27%  '
$ws.Range('E49').Value = '  -3.93%  '
$ws.Range('D50').Value = '''1.105'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').Value = '''0.06637'
$ws.Range('E51').Value = '  -2.31%  '
